# Migration: add the Vis Pesaro-Lucchese result row (with real hyperlink)
# above the existing match rows, and widen the columns so the long
# YouTube links are readable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 2 before the existing data, pushing the old rows down.
$ws.Rows.Item(2).Insert()

# Fill in the new match (order mirrors the original authoring order so
# shared-string ids line up: link, team/score, then the rest).
$ws.Cells.Item(2, 5).Value = "https://www.youtube.com/watch?v=yIB87ew2T5E&list=PLD64-55Vi5w75GRKJZRxuA33GOm4IaLDb&index=1"
$ws.Cells.Item(2, 1).Value = "Vis Pesaro-Lucchese"
$ws.Cells.Item(2, 4).Value = "0-1"
$ws.Cells.Item(2, 2).Value = 20
$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(2, 6).Value = 10
$ws.Cells.Item(2, 7).Value = 21
$ws.Cells.Item(2, 8).Value = "Destro"
$ws.Cells.Item(2, 9).Value = "vispesaro"

# Turn the Link cell into a real (clickable) hyperlink.
$ws.Hyperlinks.Add($ws.Range("E2"), "https://www.youtube.com/watch?v=yIB87ew2T5E&list=PLD64-55Vi5w75GRKJZRxuA33GOm4IaLDb&index=1") | Out-Null

# Widen the columns - A:D and F:I narrower, E (Link) very wide so the
# full YouTube URLs are visible.
$ws.Range("A1:D1").EntireColumn.ColumnWidth = 19.84
$ws.Range("E1").EntireColumn.ColumnWidth = 90.34
$ws.Range("F1:I1").EntireColumn.ColumnWidth = 19.84

# Restore the selection to where the author ended up after the edit.
$ws.Range("E10").Select() | Out-Null
